# "Contingency case 1 and 2 added"
#
# ContingencyCase1 gains the "10-11" line (it only had the "11-12" line
# before) and ContingencyCase2 loses its duplicate "10-11" line (it already
# has it at row 17, but also had it duplicated with an extra "11-12" row
# at the end) - net effect: a row moves from the end of ContingencyCase2's
# table into a new row appended to the end of ContingencyCase1's table.

$wb = $excel.ActiveWorkbook

$case1 = $wb.Worksheets.Item("ContingencyCase1")
$case2 = $wb.Worksheets.Item("ContingencyCase2")

# --- ContingencyCase1: row 17 becomes the 10-11 line, and a new row 18
#     (the 11-12 line) is appended ------------------------------------------
$case1.Cells.Item(17, 1).Value = 10
$case1.Cells.Item(17, 2).Value = 11
$case1.Cells.Item(17, 3).Value = 0.22092000000000001
$case1.Cells.Item(17, 4).Value = 0.19988
$case1.Cells.Item(17, 5).Value = 0
$case1.Cells.Item(17, 6).Value = 9999

$case1.Cells.Item(18, 1).Value = 11
$case1.Cells.Item(18, 2).Value = 12
$case1.Cells.Item(18, 3).Value = 0.17093
$case1.Cells.Item(18, 4).Value = 0.34802
$case1.Cells.Item(18, 5).Value = 0
$case1.Cells.Item(18, 6).Value = 9999

# --- ContingencyCase2: row 17 becomes the 11-12 line, and the old row 18
#     is removed -------------------------------------------------------------
$case2.Cells.Item(17, 1).Value = 11
$case2.Cells.Item(17, 2).Value = 12
$case2.Cells.Item(17, 3).Value = 0.17093
$case2.Cells.Item(17, 4).Value = 0.34802
$case2.Cells.Item(17, 5).Value = 0
$case2.Cells.Item(17, 6).Value = 9999

$case2.Rows.Item(18).Delete()

# --- Selections left behind on each sheet by the edit ----------------------
$ws1 = $wb.Worksheets.Item("BusData")
$ws2 = $wb.Worksheets.Item("LineData")

[void]$ws1.Range("G25").Select()
[void]$ws2.Range("F20").Select()
[void]$case1.Range("B25").Select()
[void]$case2.Range("F19").Select()

# ContingencyCase2 ends up the active sheet/tab after the edit.
$case2.Activate()
